# Version Control history table: the "2.8.1 / แก้ไข" row lists the
# reviewer ("ผู้ตรวจ") and responsible ("ผู้รับผิดชอบ") people for that
# revision. Update them:
#   วรรัตน์ (QM)   -> ณัฐดนัย (DM)
#   กิตติพศ (SP)   -> วิรัตน์ (TL)
#
# Both names are duplicated elsewhere in the same table (the 1.5.1 /
# จัดทำ row uses the same two names) so we must scope every Find/Replace
# to the exact table cell instead of running it over the whole story.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the row whose version column reads "2.8.1" and whose action
# column reads "แก้ไข" (there's also a "2.8.1" row higher up in the
# table header block, so the action-column check disambiguates it).
$targetRow = $null
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    if ($t.Rows.Item($r).Cells.Count -eq 5) {
        $verText = $t.Cell($r, 1).Range.Text.TrimEnd([char]7, [char]13, [char]10)
        $actionText = $t.Cell($r, 3).Range.Text.TrimEnd([char]7, [char]13, [char]10)
        if ($verText -eq "2.8.1" -and $actionText -eq "แก้ไข") {
            $targetRow = $r
        }
    }
}

# Cell 4 = ผู้รับผิดชอบ (responsible), Cell 5 = ผู้ตรวจ (reviewer)
$cellResponsible = $t.Cell($targetRow, 4)
$cellReviewer = $t.Cell($targetRow, 5)

# Re-wrap each Cell's Range through Document.Range(start, end) — using
# the Cell.Range object directly with Find/Execute is not reliably
# scoped, but a plain Range built from its Start/End is.
$rResp1 = $d.Range($cellResponsible.Range.Start, $cellResponsible.Range.End)
$rResp1.Find.Execute("วรรัตน์", $false, $false, $false, $false, $false, $true, 1, $false, "ณัฐดนัย", 2) | Out-Null

$rResp2 = $d.Range($cellResponsible.Range.Start, $cellResponsible.Range.End)
$rResp2.Find.Execute(" (QM)", $false, $false, $false, $false, $false, $true, 1, $false, " (DM)", 2) | Out-Null

$rRev1 = $d.Range($cellReviewer.Range.Start, $cellReviewer.Range.End)
$rRev1.Find.Execute("กิตติพศ ", $false, $false, $false, $false, $false, $true, 1, $false, "วิรัตน์", 2) | Out-Null

$rRev2 = $d.Range($cellReviewer.Range.Start, $cellReviewer.Range.End)
$rRev2.Find.Execute("(SP)", $false, $false, $false, $false, $false, $true, 1, $false, " (TL)", 2) | Out-Null
